$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 80019
$ws.Range("I21").Value = 80019
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 80019
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -79551
$ws.Range("N21").Value = $null
$ws.Range("H23").Value = 80019
$ws.Range("I23").Value = 80019
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 80019
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -79785
$ws.Range("N23").Value = $null
$ws.Range("H88").Value = 1820
$ws.Range("J88").Value = 1820
$ws.Range("L88").Value = 1820
$ws.Range("N88").Value = -2632
$ws.Range("H91").Value = 1820
$ws.Range("J91").Value = 1820
$ws.Range("L91").Value = 1820
$ws.Range("N91").Value = -4628
$ws.Range("H135").Value = 1017.9429
$ws.Range("I135").Value = 573
$ws.Range("J135").Value = 1489.0588
$ws.Range("K135").Value = 5157
$ws.Range("L135").Value = 13401.5292
$ws.Range("M135").Value = -2622
$ws.Range("N135").Value = -18471.5292
$ws.Range("H138").Value = 1580.5745
$ws.Range("I138").Value = 1234.7333
$ws.Range("J138").Value = 2190.8823
$ws.Range("K138").Value = 3704.199900000001
$ws.Range("L138").Value = 6572.646900000001
$ws.Range("M138").Value = 1435.800099999999
$ws.Range("N138").Value = -16852.6469

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 5000
$ws.Range("J12").Value = 5000
$ws.Range("L12").Value = 5000
$ws.Range("N12").Value = -5346
$ws.Range("H45").Value = 1839
$ws.Range("I45").Value = 1782.2222
$ws.Range("J45").Value = 2009.3334
$ws.Range("K45").Value = 1782.2222
$ws.Range("L45").Value = 2009.3334
$ws.Range("M45").Value = -1405.2222
$ws.Range("N45").Value = -2763.3334
$ws.Range("H122").Value = 1047.5264
$ws.Range("I122").Value = 841.25
$ws.Range("J122").Value = 1401.1428
$ws.Range("K122").Value = 2523.75
$ws.Range("L122").Value = 4203.428400000001
$ws.Range("M122").Value = -73.75
$ws.Range("N122").Value = -9103.428400000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 26000
$ws.Range("J112").Value = 26000
$ws.Range("L112").Value = 26000
$ws.Range("N112").Value = -28954

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 605.4375
$ws.Range("I22").Value = 583.46155
$ws.Range("K22").Value = 583.46155
$ws.Range("M22").Value = -233.46155
$ws.Range("H99").Value = 2224.6316
$ws.Range("I99").Value = 1704
$ws.Range("J99").Value = 2464.923
$ws.Range("K99").Value = 1704
$ws.Range("L99").Value = 2464.923
$ws.Range("M99").Value = -206
$ws.Range("N99").Value = -5460.923
$ws.Range("H126").Value = 2224.6316
$ws.Range("I126").Value = 1704
$ws.Range("J126").Value = 2464.923
$ws.Range("K126").Value = 5112
$ws.Range("L126").Value = 7394.768999999999
$ws.Range("M126").Value = -2642
$ws.Range("N126").Value = -12334.769

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3562.5
$ws.Range("J39").Value = 4970
$ws.Range("L39").Value = 14910
$ws.Range("N39").Value = -15498
$ws.Range("H49").Value = 3643.2856
$ws.Range("J49").Value = 3916.6667
$ws.Range("L49").Value = 11750.0001
$ws.Range("N49").Value = -12062.0001
$ws.Range("H58").Value = 1000
$ws.Range("I58").Value = 1000
$ws.Range("K58").Value = 3000
$ws.Range("M58").Value = -2872
$ws.Range("H113").Value = 739.0909
$ws.Range("I113").Value = 433.33334
$ws.Range("J113").Value = 1106
$ws.Range("K113").Value = 1300.00002
$ws.Range("L113").Value = 3318
$ws.Range("M113").Value = 869.9999800000001
$ws.Range("N113").Value = -7658

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 43753
$ws.Range("J20").Value = 52303.6
$ws.Range("L20").Value = 52303.6
$ws.Range("N20").Value = -52793.6
$ws.Range("H24").Value = 3700502.2
$ws.Range("I24").Value = 7333333.5
$ws.Range("J24").Value = 67671.336
$ws.Range("K24").Value = 7333333.5
$ws.Range("L24").Value = 67671.336
$ws.Range("M24").Value = -7333160.5
$ws.Range("N24").Value = -68017.336
$ws.Range("H102").Value = 1927.7812
$ws.Range("I102").Value = 1347.75
$ws.Range("J102").Value = 2894.5
$ws.Range("K102").Value = 1347.75
$ws.Range("L102").Value = 2894.5
$ws.Range("M102").Value = 274.25
$ws.Range("N102").Value = -6138.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3393
$ws.Range("I7").Value = 3127.2727
$ws.Range("J7").Value = 4123.75
$ws.Range("K7").Value = 3127.2727
$ws.Range("L7").Value = 4123.75
$ws.Range("M7").Value = -3015.2727
$ws.Range("N7").Value = -4347.75
$ws.Range("H21").Value = 4055.5557
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 4055.5557
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 4055.5557
$ws.Range("M21").Value = $null
$ws.Range("N21").Value = -4403.5557
$ws.Range("H22").Value = 912.8570999999999
$ws.Range("J22").Value = 981.6667
$ws.Range("L22").Value = 981.6667
$ws.Range("N22").Value = -1571.6667
$ws.Range("H27").Value = 912.8570999999999
$ws.Range("J27").Value = 981.6667
$ws.Range("L27").Value = 981.6667
$ws.Range("N27").Value = -1195.6667
$ws.Range("H34").Value = 18024
$ws.Range("J34").Value = 18024
$ws.Range("L34").Value = 18024
$ws.Range("N34").Value = -18368
$ws.Range("H40").Value = 4411.25
$ws.Range("I40").Value = 3060
$ws.Range("J40").Value = 6663.3335
$ws.Range("K40").Value = 3060
$ws.Range("L40").Value = 6663.3335
$ws.Range("M40").Value = -2924
$ws.Range("N40").Value = -6935.3335
$ws.Range("H46").Value = 1111875.5
$ws.Range("I46").Value = 622.5
$ws.Range("J46").Value = 2000878
$ws.Range("K46").Value = 622.5
$ws.Range("L46").Value = 2000878
$ws.Range("M46").Value = -434.5
$ws.Range("N46").Value = -2001254
$ws.Range("H61").Value = 4800.0713
$ws.Range("I61").Value = 4927.909
$ws.Range("J61").Value = 4331.3335
$ws.Range("K61").Value = 4927.909
$ws.Range("L61").Value = 4331.3335
$ws.Range("M61").Value = -4725.909
$ws.Range("N61").Value = -4735.3335
$ws.Range("H113").Value = 4800.0713
$ws.Range("I113").Value = 4927.909
$ws.Range("J113").Value = 4331.3335
$ws.Range("K113").Value = 4927.909
$ws.Range("L113").Value = 4331.3335
$ws.Range("M113").Value = -2757.909
$ws.Range("N113").Value = -8671.333500000001
$ws.Range("H126").Value = 3393
$ws.Range("I126").Value = 3127.2727
$ws.Range("J126").Value = 4123.75
$ws.Range("K126").Value = 9381.8181
$ws.Range("L126").Value = 12371.25
$ws.Range("M126").Value = -6911.8181
$ws.Range("N126").Value = -17311.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").Value = $null
$ws.Range("H20").Value = 46404.2
$ws.Range("I20").Value = 10999.5
$ws.Range("J20").Value = 70007.336
$ws.Range("K20").Value = 10999.5
$ws.Range("L20").Value = 70007.336
$ws.Range("M20").Value = -10759.5
$ws.Range("N20").Value = -70487.336
$ws.Range("H105").Value = 15990
$ws.Range("J105").Value = 15990
$ws.Range("L105").Value = 15990
$ws.Range("N105").Value = -22978
$ws.Range("H122").Value = 2476.8667
$ws.Range("I122").Value = 1709
$ws.Range("J122").Value = 3354.4285
$ws.Range("K122").Value = 5127
$ws.Range("L122").Value = 10063.2855
$ws.Range("M122").Value = -2677
$ws.Range("N122").Value = -14963.2855
$ws.Range("H126").Value = 3980.0667
$ws.Range("I126").Value = 4275.25
$ws.Range("J126").Value = 2799.3333
$ws.Range("K126").Value = 12825.75
$ws.Range("L126").Value = 8397.999899999999
$ws.Range("M126").Value = -10355.75
$ws.Range("N126").Value = -13337.9999
$ws.Range("H132").Value = 1766.9434
$ws.Range("I132").Value = 1428.129
$ws.Range("J132").Value = 2244.3635
$ws.Range("K132").Value = 4284.387
$ws.Range("L132").Value = 6733.0905
$ws.Range("M132").Value = -1754.387
$ws.Range("N132").Value = -11793.0905
$ws.Range("H136").Value = 1236790.8
$ws.Range("I136").Value = 3088443.5
$ws.Range("K136").Value = 9265330.5
$ws.Range("M136").Value = -9262780.5
